# Insert a new data row at row 506 (pushing existing rows 506.. down by one)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 506:540 down to 507:541, inheriting formatting from the row above.
$ws.Rows.Item(506).Insert()

# Populate the newly inserted (now blank) row 506 with the new record.
$ws.Cells.Item(506, 1).Value = 5
$ws.Cells.Item(506, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(506, 3).Value = "Maule"
$ws.Cells.Item(506, 4).Value = 45267
$ws.Cells.Item(506, 5).Value = 7
$ws.Cells.Item(506, 6).Value = 100112008
$ws.Cells.Item(506, 7).Value = "Coliflor"
$ws.Cells.Item(506, 8).Value = "Sin especificar"
$ws.Cells.Item(506, 9).Value = "Primera"
$ws.Cells.Item(506, 10).Value = 5000
$ws.Cells.Item(506, 11).Value = 1000
$ws.Cells.Item(506, 12).Value = 1000
$ws.Cells.Item(506, 13).Value = 1000
$ws.Cells.Item(506, 14).Value = "$/unidad"
$ws.Cells.Item(506, 15).Value = "Región del Maule"
$ws.Cells.Item(506, 16).Value = 1000
$ws.Cells.Item(506, 17).Value = 1
$ws.Cells.Item(506, 18).Value = "Hortaliza"
